# Common: Some experiments with graphs
# Append new translation rows (815-831) to the "Import" sheet, mirroring
# the existing key/value translation rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Carry the same row formatting (style) used by the last existing data
# row (814, wrap-text style) down through the new rows.
$ws.Range("A814:C814").Copy() | Out-Null
$ws.Range("A815:C831").PasteSpecial(-4122) | Out-Null

# Rows 815-818: language keys then labels then Czech captions
# (min/max/average/median plot labels).
$ws.Cells.Item(815, 1).Value = "cs"
$ws.Cells.Item(816, 1).Value = "cs"
$ws.Cells.Item(817, 1).Value = "cs"
$ws.Cells.Item(818, 1).Value = "cs"

$ws.Cells.Item(815, 2).Value = "lab.vape.plot.min.label"
$ws.Cells.Item(816, 2).Value = "lab.vape.plot.max.label"
$ws.Cells.Item(817, 2).Value = "lab.vape.plot.average.label"
$ws.Cells.Item(818, 2).Value = "lab.vape.plot.median.label"

$ws.Cells.Item(815, 3).Value = "Nejhorší"
$ws.Cells.Item(816, 3).Value = "Nejlepší"
$ws.Cells.Item(817, 3).Value = "Průměr"
$ws.Cells.Item(818, 3).Value = "Medián"

# Row 819: rating count label.
$ws.Cells.Item(819, 1).Value = "cs"
$ws.Cells.Item(819, 2).Value = "lab.vape.plot.count.label"
$ws.Cells.Item(819, 3).Value = "Počet hodnocení"

# Rows 820-829: plot column labels, reusing existing Czech captions.
$ws.Cells.Item(820, 1).Value = "cs"
$ws.Cells.Item(820, 2).Value = "lab.vape.plot.rating.column"
$ws.Cells.Item(820, 3).Value = "Celkové hodnocení"

$ws.Cells.Item(821, 1).Value = "cs"
$ws.Cells.Item(821, 2).Value = "lab.vape.plot.taste.column"
$ws.Cells.Item(821, 3).Value = "Chuť"

$ws.Cells.Item(822, 1).Value = "cs"
$ws.Cells.Item(822, 2).Value = "lab.vape.plot.fruits.column"
$ws.Cells.Item(822, 3).Value = "Ovocné tóny"

$ws.Cells.Item(823, 1).Value = "cs"
$ws.Cells.Item(823, 2).Value = "lab.vape.plot.complex.column"
$ws.Cells.Item(823, 3).Value = "Komplexní"

$ws.Cells.Item(824, 1).Value = "cs"
$ws.Cells.Item(824, 2).Value = "lab.vape.plot.tobacco.column"
$ws.Cells.Item(824, 3).Value = "Tabák"

$ws.Cells.Item(825, 1).Value = "cs"
$ws.Cells.Item(825, 2).Value = "lab.vape.plot.fresh.column"
$ws.Cells.Item(825, 3).Value = "Větrnost"

$ws.Cells.Item(826, 1).Value = "cs"
$ws.Cells.Item(826, 2).Value = "lab.vape.plot.mtl.column"
$ws.Cells.Item(826, 3).Value = "MTL"

$ws.Cells.Item(827, 1).Value = "cs"
$ws.Cells.Item(827, 2).Value = "lab.vape.plot.dl.column"
$ws.Cells.Item(827, 3).Value = "DL"

$ws.Cells.Item(828, 1).Value = "cs"
$ws.Cells.Item(828, 2).Value = "lab.vape.plot.throathit.column"
$ws.Cells.Item(828, 3).Value = "Throat hit"

$ws.Cells.Item(829, 1).Value = "cs"
$ws.Cells.Item(829, 2).Value = "lab.vape.plot.cakes.column"
$ws.Cells.Item(829, 3).Value = "Buchty"

# Row 830: clouds column label (brand new Czech caption "Mraky").
$ws.Cells.Item(830, 1).Value = "cs"
$ws.Cells.Item(830, 2).Value = "lab.vape.plot.clouds.column"
$ws.Cells.Item(830, 3).Value = "Mraky"

# Row 831: overall plot title.
$ws.Cells.Item(831, 1).Value = "cs"
$ws.Cells.Item(831, 2).Value = "lab.vape.plot.title"
$ws.Cells.Item(831, 3).Value = "Přehled hodnocení vapování"

# Match the saved view state: scrolled near the bottom of the sheet with
# the newly added title row selected.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 808
$win.ScrollColumn = 1
$ws.Range("B825").Select() | Out-Null
